# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers in AD, AE, AF, styled like the
# existing headers (bold, centered horizontally/vertically, thin box border)
# to match the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1        # xlContinuous (thin box border)

# Data rows 2-49 - every row gets the same team record values.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 83  # AD
    $ws.Cells.Item($r, 31).Value = 79  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
